$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.87 = 41514.46 pesos`n✅ 41514.46 pesos = 9.81 = 951.09 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the Binance transfi rate table ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 101.29
$wsTasas.Range("O10").Value = 4205
$wsTasas.Range("N12").Value = 4233.99
$wsTasas.Range("O12").Value = 97
